$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25, columns C-O (per commit diff)
$data = @{
    2 = @{ 3 = 1.000622686334603; 4 = 1.00042768552868; 5 = 1.000137389058627; 6 = 0.9999999999999999; 7 = 0.9994409678554768; 8 = 1.001335142011575; 9 = 1; 10 = 1; 11 = 0.9999999999999999; 12 = 0.9929700189459452; 13 = 0.9809305430599641; 14 = 0.9892571214978071; 15 = 0.9917783673874212 }
    3 = @{ 3 = 1.000648887587357; 4 = 1.000451453704671; 5 = 1.00016349331431; 6 = 1; 7 = 0.9995328624433907; 8 = 1.00135893174693; 9 = 0.9999999999999998; 10 = 1; 11 = 1; 12 = 0.9930806770244589; 13 = 0.9810412564835645; 14 = 0.9897355225267479; 15 = 0.9919554724000207 }
    4 = @{ 3 = 1.000673233931444; 4 = 1.000471764099104; 5 = 1.000184803473164; 6 = 1; 7 = 0.9996069307097758; 8 = 1.001379260564283; 9 = 1; 10 = 1; 11 = 1; 12 = 0.9931704736126865; 13 = 0.9811289238264292; 14 = 0.9901015465064754; 15 = 0.9920955465680324 }
    5 = @{ 3 = 1.000685687303153; 4 = 1.000482319909576; 5 = 1.000196002764639; 6 = 0.9999999999999999; 7 = 0.9996457784217335; 8 = 1.001389825949598; 9 = 1; 10 = 1; 11 = 1; 12 = 0.9932169198772401; 13 = 0.9811744545551909; 14 = 0.9902896941478881; 15 = 0.9921686740540031 }
    6 = @{ 3 = 1.000690674793759; 4 = 1.000486745589122; 5 = 1.000200823596734; 6 = 1; 7 = 0.9996624957564926; 8 = 1.00139425564354; 9 = 1; 10 = 0.9999999999999999; 11 = 0.9999999999999998; 12 = 0.9932365737097123; 13 = 0.9811939685496838; 14 = 0.9903704100706129; 15 = 0.9922001828926909 }
    7 = @{ 3 = 1.0006863141471; 4 = 1.000484981855278; 5 = 1.000199980933393; 6 = 1; 7 = 0.9996607263622661; 8 = 1.001392490309868; 9 = 0.9999999999999999; 10 = 1; 11 = 0.9999999999999999; 12 = 0.993232178275632; 13 = 0.9811921686871382; 14 = 0.9903695559363657; 15 = 0.9921983993505202 }
    8 = @{ 3 = 1.000683736405587; 4 = 1.00048335451887; 5 = 1.000198586415924; 6 = 1; 7 = 0.999656058688524; 8 = 1.001390861497355; 9 = 0.9999999999999999; 10 = 1; 11 = 1; 12 = 0.9932259697655483; 13 = 0.9811868063691432; 14 = 0.9903494526298356; 15 = 0.9921900765318403 }
    9 = @{ 3 = 1.000666468136869; 4 = 1.000467842017148; 5 = 1.000181866556255; 6 = 1; 7 = 0.9995965293664059; 8 = 1.001375334924729; 9 = 1; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9931558529481095; 13 = 0.9811169206707427; 14 = 0.9900576725936259; 15 = 0.9920772434244498 }
    10 = @{ 3 = 1.000625294634556; 4 = 1.000429804225291; 5 = 1.000139755932562; 6 = 0.9999999999999998; 7 = 0.9994481406583978; 8 = 1.001337262629988; 9 = 1; 10 = 1; 11 = 1; 12 = 0.9929789542249997; 13 = 0.9809391779520938; 14 = 0.9892936173015655; 15 = 0.9917919215754062 }
    11 = @{ 3 = 1.000582054355209; 4 = 1.000391595425291; 5 = 1.000097754307588; 6 = 0.9999999999999999; 7 = 0.9993024243357699; 8 = 1.001299019171988; 9 = 0.9999999999999999; 10 = 0.9999999999999999; 11 = 0.9999999999999998; 12 = 0.9927993114579459; 13 = 0.9807604900997736; 14 = 0.9885064957960398; 15 = 0.9915085623489757 }
    12 = @{ 3 = 1.000559426639291; 4 = 1.000372491282045; 5 = 1.00007705577667; 6 = 0.9999999999999998; 7 = 0.9992315333092147; 8 = 1.001279897699974; 9 = 1; 10 = 1; 11 = 0.9999999999999999; 12 = 0.99270960401668; 13 = 0.9806722775622932; 14 = 0.9881130276643748; 15 = 0.9913699747436218 }
    13 = @{ 3 = 1.000541244143242; 4 = 1.000357471265527; 5 = 1.000060753421116; 6 = 1; 7 = 0.9991760424823473; 8 = 1.001264864059272; 9 = 0.9999999999999998; 10 = 0.9999999999999997; 11 = 1; 12 = 0.9926384098310759; 13 = 0.9806026348348821; 14 = 0.9877994012687674; 15 = 0.9912609805939937 }
    14 = @{ 3 = 1.000537865757148; 4 = 1.000354864844444; 5 = 1.000057949301194; 6 = 0.9999999999999998; 7 = 0.999166323554609; 8 = 1.001262255273986; 9 = 0.9999999999999999; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9926256613463961; 13 = 0.9805903750571657; 14 = 0.9877438106887179; 15 = 0.9912418059023099 }
    15 = @{ 3 = 1.000542284819078; 4 = 1.000358721834498; 5 = 1.000062225962589; 6 = 1; 7 = 0.9991804306813162; 8 = 1.001266115762595; 9 = 0.9999999999999999; 10 = 1; 11 = 0.9999999999999999; 12 = 0.9926435746617879; 13 = 0.9806081404546154; 14 = 0.9878241117195113; 15 = 0.9912695355450851 }
    16 = @{ 3 = 1.000542092452296; 4 = 1.000358612251927; 5 = 1.000062063524042; 6 = 0.9999999999999999; 7 = 0.9991802229992125; 8 = 1.001266006080625; 9 = 1; 10 = 0.9999999999999999; 11 = 0.9999999999999998; 12 = 0.9926432523169973; 13 = 0.9806078967258739; 14 = 0.987823222842206; 15 = 0.9912691973514873 }
    17 = @{ 3 = 1.000541652408659; 4 = 1.000358286554168; 5 = 1.000061650701614; 6 = 0.9999999999999999; 7 = 0.9991792725455588; 8 = 1.001265680087436; 9 = 0.9999999999999999; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9926419896014912; 13 = 0.9806067227458315; 14 = 0.9878181849707141; 15 = 0.9912674171440682 }
    18 = @{ 3 = 1.000555976898579; 4 = 1.000370509319813; 5 = 1.000075032589469; 6 = 0.9999999999999999; 7 = 0.999224614024952; 8 = 1.001277913939967; 9 = 0.9999999999999998; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9926997039777591; 13 = 0.980663657766065; 14 = 0.9880750142808115; 15 = 0.9913565553642379 }
    19 = @{ 3 = 1.000564856698452; 4 = 1.000378000601609; 5 = 1.000083169866812; 6 = 0.9999999999999999; 7 = 0.999252533221669; 8 = 1.001285412016869; 9 = 0.9999999999999998; 10 = 1; 11 = 1; 12 = 0.9927351777565511; 13 = 0.9806985490841018; 14 = 0.9882315646229769; 15 = 0.9914113147767259 }
    20 = @{ 3 = 1.000565649895043; 4 = 1.000378103869158; 5 = 1.000083159769493; 6 = 0.9999999999999999; 7 = 0.9992524866077929; 8 = 1.001285515378089; 9 = 1; 10 = 1; 11 = 1; 12 = 0.992735783090937; 13 = 0.9806984544175619; 14 = 0.9882304672405625; 15 = 0.9914110723598614 }
    21 = @{ 3 = 1.00054239099987; 4 = 1.000358201871338; 5 = 1.000061505170543; 6 = 0.9999999999999998; 7 = 0.9991791271971745; 8 = 1.001265595327792; 9 = 1; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9926426555460733; 13 = 0.9806065548145206; 14 = 0.9878175896396908; 15 = 0.9912671910052615 }
    22 = @{ 3 = 1.000518484892692; 4 = 1.000338356618687; 5 = 1.000039924408681; 6 = 1; 7 = 0.9991068376800544; 8 = 1.001245732074136; 9 = 0.9999999999999998; 10 = 1; 11 = 1; 12 = 0.9925493520037492; 13 = 0.980515185785389; 14 = 0.9874034852814022; 15 = 0.9911248500884366 }
    23 = @{ 3 = 1.000499703369272; 4 = 1.000322888296105; 5 = 1.000023061129193; 6 = 0.9999999999999999; 7 = 0.9990508584695077; 8 = 1.001230249720727; 9 = 1; 10 = 0.9999999999999999; 11 = 0.9999999999999999; 12 = 0.9924764089133986; 13 = 0.9804438867961116; 14 = 0.9870776704419357; 15 = 0.9910141921551594 }
    24 = @{ 3 = 1.000492637765728; 4 = 1.000317066453466; 5 = 1.000016668385421; 6 = 1; 7 = 0.9990298948186871; 8 = 1.001224422597276; 9 = 0.9999999999999999; 10 = 1; 11 = 0.9999999999999999; 12 = 0.9924489784594769; 13 = 0.9804170700242427; 14 = 0.9869545402507429; 15 = 0.9909726673358247 }
    25 = @{ 3 = 1.000528998179837; 4 = 1.000346929843557; 5 = 1.000048960622729; 6 = 0.9999999999999999; 7 = 0.9991385124171362; 8 = 1.001254313075509; 9 = 1; 10 = 0.9999999999999999; 11 = 1; 12 = 0.9925905363628479; 13 = 0.9805553681104977; 14 = 0.9875864685552728; 15 = 0.9911874850353689 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $data[$row][$col]
    }
}
